# Apply metadata updates to the "Metadata" worksheet (sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Step 1: shift the tail block (old rows 12-15: Description, Purpose,
# Copyright, Immutable) down by one row to make room for the new
# "Jurisdiction" row. Work from the bottom up so sources aren't clobbered
# before they're read, and copy formatting (style) along with the values.

$a15 = $ws.Range("A15").Value2
$b15 = $ws.Range("B15").Value2
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A16").Value2 = $a15
$ws.Range("B16").Value2 = $b15

$a14 = $ws.Range("A14").Value2
$b14 = $ws.Range("B14").Value2
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$ws.Range("A15").Value2 = $a14
$ws.Range("B15").Value2 = $b14

$a13 = $ws.Range("A13").Value2
$b13 = $ws.Range("B13").Value2
$ws.Range("A13:B13").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)
$ws.Range("A14").Value2 = $a13
$ws.Range("B14").Value2 = $b13

$a12 = $ws.Range("A12").Value2
$b12 = $ws.Range("B12").Value2
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)
$ws.Range("A13").Value2 = $a12
$ws.Range("B13").Value2 = $b12

# --- Step 2: write the new "Jurisdiction" row into row 12, using the same
# formatting as the neighboring body rows.
$ws.Range("A10:B10").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A12").Value2 = "Jurisdiction"
$ws.Range("B12").Value2 = ""

# --- Step 3: row 11 keeps "Contact" in column A; update column B with the
# second contact's display string (previously a duplicate of row 10).
$ws.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# --- Step 4: update the remaining in-place metadata values.
$ws.Range("B3").Value2 = "0.1.7"
$ws.Range("B6").Value2 = "draft"
$ws.Range("B8").Value2 = "2024-08-27T12:23:18-05:00"
$ws.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

$wb.Save()
